$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Prototype"
$ws.Range("B1").Value = "ID"

$ws.Range("A2").Value = "void init();"
$ws.Range("B2").Value = "IDX1"

$ws.Range("A3").Value = "int sum (int x, int y);"
$ws.Range("B3").Value = "IDX2"

$ws.Range("A4").Value = "int sub (int x, int y);"
$ws.Range("B4").Value = "IDX3"

$ws.Range("A5").Value = "int mul (int x, int y);"
$ws.Range("B5").Value = "IDX4"

$ws.Range("A6").Value = "float div(float x, float y);"
$ws.Range("B6").Value = "IDX5"
